$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2242990654205607
$ws.Range("C2").Value = 0.5015576323987538
$ws.Range("J2").Value = 0.01557632398753894
$ws.Range("O2").Value = 0.003115264797507788
$ws.Range("P2").Value = 0.161993769470405
$ws.Range("S2").Value = 0.09345794392523364
$ws.Range("B3").Value = 0.005882352941176471
$ws.Range("C3").Value = 0.03529411764705882
$ws.Range("J3").Value = 0.05882352941176471
$ws.Range("P3").Value = 0.6647058823529411
$ws.Range("S3").Value = 0.2352941176470588
$ws.Range("B6").Value = 0.05263157894736842
$ws.Range("D6").Value = 0.003759398496240601
$ws.Range("E6").Value = 0.003759398496240601
$ws.Range("F6").Value = 0.07894736842105263
$ws.Range("J6").Value = 0.2857142857142857
$ws.Range("O6").Value = 0.03759398496240601
$ws.Range("Q6").Value = 0.1203007518796992
$ws.Range("R6").Value = 0.05263157894736842
$ws.Range("S6").Value = 0.3646616541353384
$ws.Range("B7").Value = 0.1349206349206349
$ws.Range("D7").Value = 0.0119047619047619
$ws.Range("F7").Value = 0.05952380952380952
$ws.Range("J7").Value = 0.1031746031746032
$ws.Range("O7").Value = 0.01587301587301587
$ws.Range("Q7").Value = 0.2023809523809524
$ws.Range("R7").Value = 0.07936507936507936
$ws.Range("S7").Value = 0.3928571428571428
$ws.Range("B8").Value = 0.1004016064257028
$ws.Range("D8").Value = 0.01606425702811245
$ws.Range("F8").Value = 0.08433734939759036
$ws.Range("J8").Value = 0.09236947791164658
$ws.Range("O8").Value = 0.03413654618473896
$ws.Range("Q8").Value = 0.1526104417670683
$ws.Range("R8").Value = 0.1164658634538153
$ws.Range("S8").Value = 0.4036144578313253
$ws.Range("B9").Value = 0.1428571428571428
$ws.Range("D9").Value = 0.005291005291005291
$ws.Range("F9").Value = 0.04232804232804233
$ws.Range("J9").Value = 0.1005291005291005
$ws.Range("O9").Value = 0.02645502645502645
$ws.Range("Q9").Value = 0.1587301587301587
$ws.Range("R9").Value = 0.07407407407407407
$ws.Range("S9").Value = 0.4497354497354497
$ws.Range("B10").Value = 0.0993322203672788
$ws.Range("D10").Value = 0.01419031719532554
$ws.Range("E10").Value = 0.0008347245409015025
$ws.Range("F10").Value = 0.07512520868113523
$ws.Range("J10").Value = 0.1185308848080134
$ws.Range("O10").Value = 0.03171953255425709
$ws.Range("Q10").Value = 0.1853088480801336
$ws.Range("R10").Value = 0.08931552587646077
$ws.Range("S10").Value = 0.3856427378964942
$ws.Range("G11").Value = 0.14004914004914
$ws.Range("J11").Value = 0.09828009828009827
$ws.Range("K11").Value = 0.2162162162162162
$ws.Range("L11").Value = 0.5356265356265356
$ws.Range("S11").Value = 0.009828009828009828
$ws.Range("G12").Value = 0.7379912663755459
$ws.Range("J12").Value = 0.1965065502183406
$ws.Range("K12").Value = 0.01310043668122271
$ws.Range("L12").Value = 0.03930131004366812
$ws.Range("S12").Value = 0.01310043668122271
$ws.Range("G13").Value = 0.625
$ws.Range("J13").Value = 0.3214285714285715
$ws.Range("S13").Value = 0.05357142857142857
$ws.Range("F15").Value = 0.02857142857142857
$ws.Range("H15").Value = 0.1821428571428571
$ws.Range("I15").Value = 0.04285714285714286
$ws.Range("J15").Value = 0.275
$ws.Range("K15").Value = 0.07142857142857142
$ws.Range("M15").Value = 0.01785714285714286
$ws.Range("N15").Value = 0.003571428571428571
$ws.Range("O15").Value = 0.06785714285714285
$ws.Range("S15").Value = 0.3107142857142857
$ws.Range("F16").Value = 0.01666666666666667
$ws.Range("H16").Value = 0.1888888888888889
$ws.Range("I16").Value = 0.07777777777777778
$ws.Range("J16").Value = 0.35
$ws.Range("K16").Value = 0.15
$ws.Range("M16").Value = 0.02777777777777778
$ws.Range("O16").Value = 0.06666666666666667
$ws.Range("S16").Value = 0.1222222222222222
$ws.Range("F17").Value = 0.02475247524752475
$ws.Range("H17").Value = 0.198019801980198
$ws.Range("I17").Value = 0.08168316831683169
$ws.Range("J17").Value = 0.3391089108910891
$ws.Range("K17").Value = 0.1633663366336634
$ws.Range("M17").Value = 0.01485148514851485
$ws.Range("N17").Value = 0.002475247524752475
$ws.Range("O17").Value = 0.07673267326732673
$ws.Range("S17").Value = 0.09900990099009901
$ws.Range("F18").Value = 0.01421800947867299
$ws.Range("H18").Value = 0.2085308056872038
$ws.Range("I18").Value = 0.1090047393364929
$ws.Range("J18").Value = 0.3317535545023697
$ws.Range("K18").Value = 0.0995260663507109
$ws.Range("M18").Value = 0.03317535545023697
$ws.Range("O18").Value = 0.06161137440758294
$ws.Range("S18").Value = 0.1421800947867299
$ws.Range("F19").Value = 0.02844311377245509
$ws.Range("H19").Value = 0.217814371257485
$ws.Range("I19").Value = 0.08008982035928144
$ws.Range("J19").Value = 0.3338323353293413
$ws.Range("K19").Value = 0.1324850299401198
$ws.Range("M19").Value = 0.02470059880239521
$ws.Range("N19").Value = 0.0007485029940119761
$ws.Range("O19").Value = 0.07934131736526946
$ws.Range("S19").Value = 0.1025449101796407
